$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.2142857142857143
$ws.Range("C2").Value = 0.5248447204968945
$ws.Range("J2").Value = 0.02173913043478261
$ws.Range("P2").Value = 0.1459627329192547
$ws.Range("S2").Value = 0.09316770186335403

# Row 3
$ws.Range("B3").Value = 0.01666666666666667
$ws.Range("C3").Value = 0.05
$ws.Range("J3").Value = 0.05
$ws.Range("P3").Value = 0.6833333333333333

# Row 4
$ws.Range("J4").Value = 0.05172413793103448
$ws.Range("P4").Value = 0.7413793103448276
$ws.Range("S4").Value = 0.2068965517241379

# Row 6
$ws.Range("B6").Value = 0.05627705627705628
$ws.Range("D6").Value = 0.004329004329004329
$ws.Range("J6").Value = 0.3852813852813853
$ws.Range("O6").Value = 0.004329004329004329
$ws.Range("Q6").Value = 0.1471861471861472
$ws.Range("R6").Value = 0.06060606060606061
$ws.Range("S6").Value = 0.2510822510822511

# Row 7
$ws.Range("B7").Value = 0.1206896551724138
$ws.Range("D7").Value = 0.04597701149425287
$ws.Range("F7").Value = 0.05747126436781609
$ws.Range("J7").Value = 0.1551724137931035
$ws.Range("O7").Value = 0.01149425287356322
$ws.Range("Q7").Value = 0.1609195402298851
$ws.Range("R7").Value = 0.07471264367816093
$ws.Range("S7").Value = 0.3735632183908046

# Row 8
$ws.Range("B8").Value = 0.1226415094339623
$ws.Range("D8").Value = 0.01886792452830189
$ws.Range("F8").Value = 0.08176100628930817
$ws.Range("J8").Value = 0.09748427672955975
$ws.Range("O8").Value = 0.01257861635220126
$ws.Range("Q8").Value = 0.1949685534591195
$ws.Range("R8").Value = 0.05974842767295598
$ws.Range("S8").Value = 0.4119496855345912

# Row 9
$ws.Range("B9").Value = 0.1191709844559585
$ws.Range("D9").Value = 0.0155440414507772
$ws.Range("F9").Value = 0.05699481865284974
$ws.Range("J9").Value = 0.1243523316062176
$ws.Range("O9").Value = 0.02590673575129534
$ws.Range("Q9").Value = 0.2020725388601036
$ws.Range("R9").Value = 0.07253886010362694
$ws.Range("S9").Value = 0.383419689119171

# Row 10
$ws.Range("B10").Value = 0.1114391143911439
$ws.Range("D10").Value = 0.03173431734317343
$ws.Range("E10").Value = 0.001476014760147601
$ws.Range("F10").Value = 0.05239852398523985
$ws.Range("J10").Value = 0.1195571955719557
$ws.Range("O10").Value = 0.02066420664206642
$ws.Range("Q10").Value = 0.2501845018450184
$ws.Range("R10").Value = 0.07970479704797048
$ws.Range("S10").Value = 0.3328413284132841

# Row 11
$ws.Range("G11").Value = 0.1573426573426573
$ws.Range("J11").Value = 0.09090909090909091
$ws.Range("K11").Value = 0.1958041958041958
$ws.Range("L11").Value = 0.5454545454545454
$ws.Range("S11").Value = 0.01048951048951049

# Row 12
$ws.Range("G12").Value = 0.6829268292682927
$ws.Range("J12").Value = 0.2195121951219512
$ws.Range("K12").Value = 0.01219512195121951
$ws.Range("L12").Value = 0.0426829268292683
$ws.Range("S12").Value = 0.0426829268292683

# Row 13
$ws.Range("G13").Value = 0.7391304347826086
$ws.Range("J13").Value = 0.2608695652173913

# Row 15
$ws.Range("F15").Value = 0.008658008658008658
$ws.Range("H15").Value = 0.1168831168831169
$ws.Range("I15").Value = 0.06926406926406926
$ws.Range("J15").Value = 0.4155844155844156
$ws.Range("K15").Value = 0.06060606060606061
$ws.Range("O15").Value = 0.06926406926406926
$ws.Range("S15").Value = 0.2597402597402597

# Row 16
$ws.Range("F16").Value = 0.04411764705882353
$ws.Range("H16").Value = 0.1274509803921569
$ws.Range("I16").Value = 0.08333333333333333
$ws.Range("J16").Value = 0.446078431372549
$ws.Range("K16").Value = 0.08333333333333333
$ws.Range("M16").Value = 0.009803921568627451
$ws.Range("N16").Value = 0.009803921568627451
$ws.Range("O16").Value = 0.09313725490196079
$ws.Range("S16").Value = 0.1029411764705882

# Row 17
$ws.Range("F17").Value = 0.04536489151873768
$ws.Range("H17").Value = 0.1124260355029586
$ws.Range("I17").Value = 0.09861932938856016
$ws.Range("J17").Value = 0.4536489151873767
$ws.Range("K17").Value = 0.106508875739645
$ws.Range("M17").Value = 0.01775147928994083
$ws.Range("O17").Value = 0.0631163708086785
$ws.Range("S17").Value = 0.1025641025641026

# Row 18
$ws.Range("F18").Value = 0.02366863905325444
$ws.Range("H18").Value = 0.1597633136094675
$ws.Range("I18").Value = 0.106508875739645
$ws.Range("J18").Value = 0.4556213017751479
$ws.Range("K18").Value = 0.1183431952662722
$ws.Range("O18").Value = 0.05325443786982249
$ws.Range("S18").Value = 0.08284023668639054

# Row 19
$ws.Range("F19").Value = 0.02353966870095902
$ws.Range("H19").Value = 0.1612903225806452
$ws.Range("I19").Value = 0.08195292066259809
$ws.Range("J19").Value = 0.3958151700087184
$ws.Range("K19").Value = 0.1054925893635571
$ws.Range("M19").Value = 0.01220575414123801
$ws.Range("N19").Value = 0.0008718395815170009
$ws.Range("O19").Value = 0.07933740191804708
$ws.Range("S19").Value = 0.1394943330427201
